$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# --- Insert a new row at 59 (pushes old rows 59-72 down to 60-73) ---
$ws.Rows.Item(59).Insert()

# Copy formatting from the row above (row 58), which already carries the
# look (font/border/alignment) used for the other rows of this parameter
# table, so the newly inserted row matches its neighbours visually.
$ws.Range("A58:H58").Copy()
$ws.Range("A59:H59").PasteSpecial(-4122)

# --- Populate the new row 59 with the new "AutoBatchFg" parameter ---
$ws.Cells.Item(59, 1).Value = 51
$ws.Cells.Item(59, 2).Value = "AutoBatchFg"
$ws.Cells.Item(59, 3).Value = "自動批次記號"
$ws.Cells.Item(59, 4).Value = "VARCHAR2"
$ws.Cells.Item(59, 5).Value = 1
$ws.Cells.Item(59, 7).Value = "Y:啟用" + [char]10 + "N:不啟用"
$ws.Cells.Item(59, 8).Value = "2022/4/12智偉新增"

# The new row's description (col G) wraps onto two lines, same as the
# other "Y/N flag" rows in this sheet (e.g. row 52) -- match their height.
$ws.Rows.Item(59).RowHeight = 32.4

# --- Renumber column A for the rows that shifted down (they hold static
#     numbers, not formulas, so Insert() does not bump them automatically) ---
$ws.Cells.Item(60, 1).Value = 52
$ws.Cells.Item(61, 1).Value = 53
$ws.Cells.Item(62, 1).Value = 54
$ws.Cells.Item(63, 1).Value = 55

# --- Keep the selection where the author left it after the edit ---
$ws.Range("B59").Select()
